$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new customer row (row 9): phone 79174404 with 0 points, birthday unknown.
# Phone numbers in this sheet are stored as text even though they look numeric, so
# force text formatting while entering the value, then drop the formatting again so
# the cell ends up with no explicit style applied (matching the rest of the sheet).
$phoneCell = $ws.Cells.Item(9, 1)
$phoneCell.NumberFormat = "@"
$phoneCell.Value = "79174404"
$phoneCell.ClearFormats()

# Birthday is unknown for this customer -> cell stays blank, but still present in the row
# (matching the pattern used by the other customers whose birthday isn't on file).
$birthdayCell = $ws.Cells.Item(9, 2)
$birthdayCell.Font.Bold = $false
$birthdayCell.ClearFormats()

# Points total for this customer.
$ws.Cells.Item(9, 3).Value = 0
